# Add a new "UK" worksheet (Test Data for UK Market), based on the
# existing "Poland" sheet template, and make it the active sheet.

$wb = $excel.ActiveWorkbook

# The "Poland" sheet is the last market sheet and acts as the template
# for every new market tab that gets added to this gallery workbook.
$template = $wb.Worksheets.Item("Poland")

# Copy it, placing the new sheet right after the template (i.e. at the
# end of the workbook, becoming the new last tab).
$template.Copy($null, $template)

# The freshly created copy is now the last worksheet in the workbook.
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Fill in the market-specific values (mirrors how every other market
# tab is populated): User Story ticket reference and market name.
$newSheet.Range("B4").Value = "NGC-2741/T3339"
$newSheet.Range("B2").Value = "UK Market"

# Match the saved selection/active cell for the new sheet.
$newSheet.Range("B4").Select() | Out-Null
